# feat: add 2022-Q4 data
#
#  1. Insert a brand new worksheet named "2022-Q4" right after "总计" (i.e. it becomes
#     the 2nd sheet, pushing "2022-Q2" and everything after it one slot to the right).
#  2. Fill that sheet with the fund-holdings table for 2022-Q4 (same layout as the
#     other quarter sheets: header in B1:H1, data starting row 2, index in col A).
#  3. Insert a new row at the top of the "总计" (summary) sheet's data table for the
#     2022-Q4 totals, shifting the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet positioned right after "总计".
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Populate the new "2022-Q4" sheet.
# ---------------------------------------------------------------------------

# Header row (bold, centered, bordered) -- matches the sibling quarter sheets.
$newSheet.Range("B1:H1").NumberFormat = "@"
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"
$newSheet.Range("B1:H1").Font.Bold = $true
$newSheet.Range("B1:H1").HorizontalAlignment = -4108
$newSheet.Range("B1:H1").VerticalAlignment = -4160
$newSheet.Range("B1:H1").Borders.LineStyle = 1

# Data rows. Columns B (code), C (name), D..G (all the numeric-looking figures) are
# stored as *text* in the source workbook (keeps values like "001891"/"010889" from
# losing their leading zeros, and mirrors "3.26" etc. verbatim) -- column A (index)
# and H (rank) are real numbers.
$newSheet.Range("B2:G6").NumberFormat = "@"

$q4Rows = @(
    @(0, "160143", "南方创业板2年定期开放混合", "3.26", "92.67", "3.19", "0.1040", 7),
    @(1, "166020", "中欧成长优选回报灵活配置混合A", "2.73", "93.51", "2.79", "0.0762", 10),
    @(2, "164826", "工银瑞信创业板两年定期开放混合A", "1.64", "80.79", "4.29", "0.0704", 6),
    @(3, "001891", "中欧成长优选回报灵活配置混合E", "2.34", "93.51", "2.79", "0.0653", 10),
    @(4, "010889", "工银瑞信创业板两年定期开放混合C", "0.14", "80.79", "4.29", "0.0060", 6)
)

$r = 2
foreach ($row in $q4Rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Index column (A2:A6) uses the same bold/centered/bordered style as the header.
$newSheet.Range("A2:A6").Font.Bold = $true
$newSheet.Range("A2:A6").HorizontalAlignment = -4108
$newSheet.Range("A2:A6").VerticalAlignment = -4160
$newSheet.Range("A2:A6").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 3. Insert the 2022-Q4 row into the "总计" summary sheet.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Clear any formatting that Insert() may have copied down from the header row --
# only column A carries the bold/centered index style in this sheet.
$totalSheet.Range("B2:D2").Font.Bold = $false
$totalSheet.Range("B2:D2").HorizontalAlignment = -4108
$totalSheet.Range("B2:D2").Borders.LineStyle = -4142

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 5
$totalSheet.Cells.Item(2, 4).Value = 0.32

$totalSheet.Cells.Item(2, 1).Font.Bold = $true
$totalSheet.Cells.Item(2, 1).HorizontalAlignment = -4108
$totalSheet.Cells.Item(2, 1).VerticalAlignment = -4160
$totalSheet.Cells.Item(2, 1).Borders.LineStyle = 1

# Renumber the index column (A) for the rows that just shifted down, so it stays
# a plain 0-based sequence matching its row position.
for ($i = 3; $i -le 9; $i++) {
    $totalSheet.Cells.Item($i, 1).Value = $i - 2
}
